$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last two data rows (old "MuSCs" row 4 and old "Resolving-Mac"
# row 5). Excel shifts the remaining rows up and -- because nothing else in
# the sheet references the "Resolving-Mac" string once its only row is gone
# -- it drops out of the shared-strings table automatically.
$ws.Rows("4:5").Delete()

# Row 2 used to be the "ECs" target-cluster row; it is now replaced with the
# freshly recomputed "FAPs" row (re-run with the updated TPM data).
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.7349876666666667
$ws.Range("H2").Value = 2.204963
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.692918333333334
$ws.Range("N2").Value = 11.078755
$ws.Range("O2").Value = 0.9367638696585037
$ws.Range("P2").Value = 0.956934821355692
$ws.Range("Q2").Value = 2.714249429007222
$ws.Range("R2").Value = 24.428244861065
$ws.Range("S2").Value = 0.9367638696585037
$ws.Range("T2").Value = 0.956934821355692

# Row 3 used to be the "FAPs" target-cluster row; it is now replaced with the
# freshly recomputed "MuSCs" row.
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.7349876666666667
$ws.Range("H3").Value = 2.204963
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.24929
$ws.Range("N3").Value = 0.49858
$ws.Range("O3").Value = 0.06323613034149642
$ws.Range("P3").Value = 0.04306517864430803
$ws.Range("Q3").Value = 0.1832250754233334
$ws.Range("R3").Value = 1.09935045254
$ws.Range("S3").Value = 0.06323613034149642
$ws.Range("T3").Value = 0.04306517864430803
